# Natmi following Dr Hou advice
# Update App -> Lrp1 LR-pair expression/specificity metrics for rows 2-10
# (Ligand-expressing cells count changed from 1 to 3, with corresponding
# recalculated expression/specificity values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 104.794801
$ws.Range("H2").Value = 314.384403
$ws.Range("I2").Value = 0.3872421191355361
$ws.Range("J2").Value = 0.3872421191355361
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 1798.96309521053
$ws.Range("R2").Value = 16190.66785689477
$ws.Range("S2").Value = 0.02169893390317983
$ws.Range("T2").Value = 0.02169893390317983
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 104.794801
$ws.Range("H3").Value = 314.384403
$ws.Range("I3").Value = 0.3872421191355361
$ws.Range("J3").Value = 0.3872421191355361
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 26874.02994498991
$ws.Range("R3").Value = 241866.2695049091
$ws.Range("S3").Value = 0.3241521746838101
$ws.Range("T3").Value = 0.3241521746838101
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 104.794801
$ws.Range("H4").Value = 314.384403
$ws.Range("I4").Value = 0.3872421191355361
$ws.Range("J4").Value = 0.3872421191355361
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 3431.546488990989
$ws.Range("R4").Value = 30883.9184009189
$ws.Range("S4").Value = 0.04139101054854615
$ws.Range("T4").Value = 0.04139101054854615
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 136.674446
$ws.Range("H5").Value = 410.023338
$ws.Range("I5").Value = 0.5050451128841343
$ws.Range("J5").Value = 0.5050451128841343
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 2346.22597749238
$ws.Range("R5").Value = 21116.03379743142
$ws.Range("S5").Value = 0.02829997043467567
$ws.Range("T5").Value = 0.02829997043467567
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 136.674446
$ws.Range("H6").Value = 410.023338
$ws.Range("I6").Value = 0.5050451128841343
$ws.Range("J6").Value = 0.5050451128841343
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 35049.38336128817
$ws.Range("R6").Value = 315444.4502515935
$ws.Range("S6").Value = 0.4227625652402829
$ws.Range("T6").Value = 0.4227625652402829
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 136.674446
$ws.Range("H7").Value = 410.023338
$ws.Range("I7").Value = 0.5050451128841343
$ws.Range("J7").Value = 0.5050451128841343
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 4475.457855071345
$ws.Range("R7").Value = 40279.1206956421
$ws.Range("S7").Value = 0.05398257720917568
$ws.Range("T7").Value = 0.05398257720917568
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.14904533333333
$ws.Range("H8").Value = 87.447136
$ws.Range("I8").Value = 0.1077127679803296
$ws.Range("J8").Value = 0.1077127679803296
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 500.3879611860266
$ws.Range("R8").Value = 4503.49165067424
$ws.Range("S8").Value = 0.006035635374972394
$ws.Range("T8").Value = 0.006035635374972394
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.14904533333333
$ws.Range("H9").Value = 87.447136
$ws.Range("I9").Value = 0.1077127679803296
$ws.Range("J9").Value = 0.1077127679803296
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 7475.106681636508
$ws.Range("R9").Value = 67275.96013472858
$ws.Range("S9").Value = 0.09016407631478747
$ws.Range("T9").Value = 0.09016407631478747
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.14904533333333
$ws.Range("H10").Value = 87.447136
$ws.Range("I10").Value = 0.1077127679803296
$ws.Range("J10").Value = 0.1077127679803296
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 954.4968187022861
$ws.Range("R10").Value = 8590.471368320575
$ws.Range("S10").Value = 0.0115130562905697
$ws.Range("T10").Value = 0.0115130562905697
